# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Periodo Mora" values for the two detail rows (2303 <-> 2302)
$ws.Range("E16").Value = "2302"
$ws.Range("E17").Value = "2303"

# Update "Salario Basico" amounts for both detail rows
$ws.Range("G16").Value = 1000000
$ws.Range("G17").Value = 1000000
